$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 2249.5454
$ws.Range("I18").Value = 2030.625
$ws.Range("K18").Value = 2030.625
$ws.Range("M18").Value = -1746.625

# Row 40
$ws.Range("H40").Value = 2500
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3350

# Row 42
$ws.Range("H42").Value = 219.5
$ws.Range("I42").Value = 219.5
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 658.5
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -428.5

# Row 101
$ws.Range("H101").Value = 555
$ws.Range("I101").Value = 555
$ws.Range("K101").Value = 1665
$ws.Range("M101").Value = -43

# Row 107
$ws.Range("H107").Value = 425.6875
$ws.Range("I107").Value = 454.63635
$ws.Range("K107").Value = 454.63635
$ws.Range("M107").Value = 1465.36365

# Row 132
$ws.Range("H132").Value = 1346.6
$ws.Range("I132").Value = 1234.4117
$ws.Range("J132").Value = 1982.3334
$ws.Range("K132").Value = 3703.2351
$ws.Range("L132").Value = 5947.0002
$ws.Range("M132").Value = -1173.2351
$ws.Range("N132").Value = -11007.0002

# Row 138
$ws.Range("H138").Value = 5495.636
$ws.Range("I138").Value = 4900
$ws.Range("J138").Value = 5555.2
$ws.Range("K138").Value = 14700
$ws.Range("L138").Value = 16665.6
$ws.Range("M138").Value = -9560
$ws.Range("N138").Value = -26945.6

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1552.75
$ws.Range("I2").Value = 1505.5
$ws.Range("K2").Value = 1505.5
$ws.Range("M2").Value = -1392.5

# Row 61
$ws.Range("H61").Value = 1793.3636
$ws.Range("I61").Value = 1793.3636
$ws.Range("K61").Value = 1793.3636
$ws.Range("M61").Value = -1581.3636

# Row 116
$ws.Range("H116").Value = 1552.75
$ws.Range("I116").Value = 1505.5
$ws.Range("K116").Value = 1505.5
$ws.Range("M116").Value = 788.5

# Row 136
$ws.Range("H136").Value = 1793.3636
$ws.Range("I136").Value = 1793.3636
$ws.Range("K136").Value = 5380.0908
$ws.Range("M136").Value = -2830.0908

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1552.75
$ws.Range("I3").Value = 1505.5
$ws.Range("K3").Value = 1505.5
$ws.Range("M3").Value = -1391.5

# Row 99
$ws.Range("H99").Value = 1608.4286
$ws.Range("I99").Value = 1460.6666
$ws.Range("K99").Value = 1460.6666
$ws.Range("M99").Value = 37.33339999999998

$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 40000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 40000
$ws.Range("K50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("M50").Value = 40000
$ws.Range("N50").Value = -41250

# Row 51
$ws.Range("H51").Value = 39999
$ws.Range("J51").Value = 39999
$ws.Range("L51").Value = 39999
$ws.Range("N51").Value = -41471

# Row 55
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("N55").Value = 0

# Row 58
$ws.Range("H58").Value = 5964.143
$ws.Range("I58").Value = 3492.6667
$ws.Range("K58").Value = 3492.6667
$ws.Range("M58").Value = -3289.6667

# Row 61
$ws.Range("H61").Value = 39999
$ws.Range("J61").Value = 39999
$ws.Range("L61").Value = 39999
$ws.Range("N61").Value = -40695

# Row 136
$ws.Range("H136").Value = 5964.143
$ws.Range("I136").Value = 3492.6667
$ws.Range("K136").Value = 10478.0001
$ws.Range("M136").Value = -7928.000100000001

$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 801
$ws.Range("I18").Value = 801
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 2403
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -2234

# Row 26
$ws.Range("H26").Value = 173.57143
$ws.Range("I26").Value = 173.57143
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 520.71429
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -232.71429

# Row 34
$ws.Range("H34").Value = 2002.3636
$ws.Range("I34").Value = 1717.2
$ws.Range("J34").Value = 2240
$ws.Range("K34").Value = 5151.6
$ws.Range("L34").Value = 6720
$ws.Range("M34").Value = -5067.6
$ws.Range("N34").Value = -6888

# Row 86
$ws.Range("H86").Value = 377.8
$ws.Range("I86").Value = 344.5
$ws.Range("J86").Value = 400
$ws.Range("K86").Value = 1033.5
$ws.Range("L86").Value = 1200
$ws.Range("M86").Value = 152.5
$ws.Range("N86").Value = -3572

# Row 89
$ws.Range("H89").Value = 377.8
$ws.Range("I89").Value = 344.5
$ws.Range("J89").Value = 400
$ws.Range("K89").Value = 3100.5
$ws.Range("L89").Value = 3600
$ws.Range("M89").Value = 2827.5
$ws.Range("N89").Value = -15456

# Row 104
$ws.Range("H104").Value = 72477.14
$ws.Range("I104").Value = 1384
$ws.Range("J104").Value = 250210
$ws.Range("K104").Value = 4152
$ws.Range("L104").Value = 750630
$ws.Range("M104").Value = -1531
$ws.Range("N104").Value = -755872

# Row 112
$ws.Range("H112").Value = 5000
$ws.Range("J112").Value = 5000
$ws.Range("L112").Value = 15000
$ws.Range("N112").Value = -17216

# Row 132
$ws.Range("H132").Value = 4273.8
$ws.Range("I132").Value = 967.25
$ws.Range("J132").Value = 17500
$ws.Range("K132").Value = 8705.25
$ws.Range("L132").Value = 157500
$ws.Range("M132").Value = -6175.25
$ws.Range("N132").Value = -162560

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 9829.5

# Row 83
$ws.Range("H83").Value = 9829.5

# Row 102
$ws.Range("H102").Value = 1376.4546
$ws.Range("I102").Value = 746.3684
$ws.Range("K102").Value = 746.3684
$ws.Range("M102").Value = 875.6316

# Row 132
$ws.Range("H132").Value = 3021.55
$ws.Range("J132").Value = 4622.5
$ws.Range("L132").Value = 13867.5
$ws.Range("N132").Value = -18927.5

# Row 136
$ws.Range("H136").Value = 56999.6
$ws.Range("J136").Value = 56999.6
$ws.Range("L136").Value = 170998.8
$ws.Range("N136").Value = -176098.8

$ws = $wb.Worksheets.Item("LTW")
# Row 39
$ws.Range("H39").Value = 88000
$ws.Range("J39").Value = 88000
$ws.Range("L39").Value = 88000
$ws.Range("N39").Value = -88920

# Row 46
$ws.Range("H46").Value = 3309.375
$ws.Range("I46").Value = 1994.4445
$ws.Range("K46").Value = 1994.4445
$ws.Range("M46").Value = -1806.4445

# Row 61
$ws.Range("H61").Value = 6312.5
$ws.Range("I61").Value = 7083
$ws.Range("J61").Value = 4001
$ws.Range("K61").Value = 7083
$ws.Range("L61").Value = 4001
$ws.Range("M61").Value = -6881
$ws.Range("N61").Value = -4405

# Row 82
$ws.Range("H82").Value = 2887.4443
$ws.Range("I82").Value = 2826.7144
$ws.Range("J82").Value = 3100
$ws.Range("K82").Value = 2826.7144
$ws.Range("L82").Value = 3100
$ws.Range("M82").Value = -2465.7144
$ws.Range("N82").Value = -3822

# Row 85
$ws.Range("H85").Value = 2887.4443
$ws.Range("I85").Value = 2826.7144
$ws.Range("J85").Value = 3100
$ws.Range("K85").Value = 2826.7144
$ws.Range("L85").Value = 3100
$ws.Range("M85").Value = -1578.7144
$ws.Range("N85").Value = -5596

# Row 113
$ws.Range("H113").Value = 6312.5
$ws.Range("I113").Value = 7083
$ws.Range("J113").Value = 4001
$ws.Range("K113").Value = 7083
$ws.Range("L113").Value = 4001
$ws.Range("M113").Value = -4913
$ws.Range("N113").Value = -8341

# Row 132
$ws.Range("H132").Value = 2991.0322
$ws.Range("I132").Value = 1967.5834
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 5902.7502
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -3372.7502
$ws.Range("N132").Value = -24560

$ws = $wb.Worksheets.Item("WVR")
# Row 125
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840

# Row 132
$ws.Range("H132").Value = 2747.2173
$ws.Range("I132").Value = 2286.5625
$ws.Range("J132").Value = 3800.1428
$ws.Range("K132").Value = 6859.6875
$ws.Range("L132").Value = 11400.4284
$ws.Range("M132").Value = -4329.6875
$ws.Range("N132").Value = -16460.4284

# Row 136
$ws.Range("H136").Value = 24322
$ws.Range("I136").Value = 1573.7778
$ws.Range("J136").Value = 60451.53
$ws.Range("K136").Value = 4721.3334
$ws.Range("L136").Value = 181354.59
$ws.Range("M136").Value = -2171.3334
$ws.Range("N136").Value = -186454.59
